$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new "RecvStarReward" row by inserting a blank row
# before the current row 7 (Flag / State shift down by one).
$ws.Rows.Item(7).Insert()

# Row 4: Order -> TopFinishStageOrder
$ws.Range("A4").Value = "TopFinishStageOrder"
$ws.Range("B4").Value = "INT"

# Row 5 & 6 swap: LastPlayStageNum / TopFinishStageNum
$ws.Range("A5").Value = "TopFinishStageNum"
$ws.Range("B5").Value = "INT"

$ws.Range("A6").Value = "LastPlayStageNum"
$ws.Range("B6").Value = "INT"

# Row 7 (new): RecvStarReward
$ws.Range("A7").Value = "RecvStarReward"
$ws.Range("B7").Value = "INT"

# Row 8: Flag (unchanged content, already shifted down by the insert)
$ws.Range("A8").Value = "Flag"
$ws.Range("B8").Value = "BIGINT UNSIGNED"

# Row 9: State (unchanged content, already shifted down by the insert)
$ws.Range("A9").Value = "State"
$ws.Range("B9").Value = "INT"

# Update the selection to match the saved workbook state
$ws.Range("B7").Select()
